$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update assignment 7 (PA7) deadline
$ws.Range("B9").Value = "November 13, 2024"

# Move active selection to B10, matching post-edit state
$ws.Range("B10").Select()
